# Insert a new data row at row 725 (pushing the existing rows 725:815 down
# to 726:816) and populate the new row with the values below. This mirrors
# a new weekly record being added in the middle of the price history table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 725 downwards (row 815 -> 816) to make room for the new record.
$ws.Rows("725:725").Insert()

# Populate the freshly inserted row 725 with the new observation.
$ws.Range("A725").Value = 3
$ws.Range("B725").Value = "Femacal de La Calera"
$ws.Range("C725").Value = "Coquimbo"
$ws.Range("D725").Value = 45142
$ws.Range("D725").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E725").Value = 5
$ws.Range("F725").Value = 100112032
$ws.Range("G725").Value = "Zapallo italiano"
$ws.Range("H725").Value = "Sin especificar"
$ws.Range("I725").Value = "Primera"
$ws.Range("J725").Value = 115
$ws.Range("K725").Value = 14000
$ws.Range("L725").Value = 14500
$ws.Range("M725").Value = 14283
$ws.Range("N725").Value = '$/caja 60 unidades'
$ws.Range("O725").Value = "Región de Arica y Parinacota"
$ws.Range("P725").Value = 238
$ws.Range("Q725").Value = 60
$ws.Range("R725").Value = "Hortaliza"
